$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting existing rows 250-273 down to 251-274
$ws.Rows.Item(250).Insert()

# Populate the newly inserted row 250 with the new data record
$ws.Cells.Item(250, 1).Value = 10
$ws.Cells.Item(250, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(250, 3).Value = "La Araucanía"
$ws.Cells.Item(250, 4).Value = 44461
$ws.Cells.Item(250, 5).Value = 9
$ws.Cells.Item(250, 6).Value = "Fruta"
$ws.Cells.Item(250, 7).Value = 100108
$ws.Cells.Item(250, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(250, 9).Value = 100108005
$ws.Cells.Item(250, 10).Value = "Piña"
$ws.Cells.Item(250, 11).Value = "Caramelo"
$ws.Cells.Item(250, 12).Value = "Segunda"
$ws.Cells.Item(250, 13).Value = 90
$ws.Cells.Item(250, 14).Value = 22000
$ws.Cells.Item(250, 15).Value = 23000
$ws.Cells.Item(250, 16).Value = 22556
$ws.Cells.Item(250, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(250, 18).Value = "Ecuador"
$ws.Cells.Item(250, 19).Value = 1611
$ws.Cells.Item(250, 20).Value = 14
